$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell F1 - match formatting of existing header cell (E1)
$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$timeValues = @(
    "2021-10-05 10:52:48.770374",
    "2021-10-05 10:52:48.770386",
    "2021-10-05 10:52:48.770389",
    "2021-10-05 10:52:48.770392",
    "2021-10-05 10:52:48.770395",
    "2021-10-05 10:52:48.770398",
    "2021-10-05 10:52:48.770400",
    "2021-10-05 10:52:48.770403",
    "2021-10-05 10:52:48.770406",
    "2021-10-05 10:52:48.770409",
    "2021-10-05 10:52:48.770411",
    "2021-10-05 10:52:48.770414",
    "2021-10-05 10:52:48.770416",
    "2021-10-05 10:52:48.770419",
    "2021-10-05 10:52:48.770421",
    "2021-10-05 10:52:48.770424",
    "2021-10-05 10:52:48.770427",
    "2021-10-05 10:52:48.770430",
    "2021-10-05 10:52:48.770433",
    "2021-10-05 10:52:48.770435",
    "2021-10-05 10:52:48.770438",
    "2021-10-05 10:52:48.770440",
    "2021-10-05 10:52:48.770443",
    "2021-10-05 10:52:48.770445",
    "2021-10-05 10:52:48.770448",
    "2021-10-05 10:52:48.770451",
    "2021-10-05 10:52:48.770454",
    "2021-10-05 10:52:48.770456",
    "2021-10-05 10:52:48.770459",
    "2021-10-05 10:52:48.770462",
    "2021-10-05 10:52:48.770464",
    "2021-10-05 10:52:48.770467",
    "2021-10-05 10:52:48.770470",
    "2021-10-05 10:52:48.770473",
    "2021-10-05 10:52:48.770475",
    "2021-10-05 10:52:48.770478",
    "2021-10-05 10:52:48.770481",
    "2021-10-05 10:52:48.770484",
    "2021-10-05 10:52:48.770486",
    "2021-10-05 10:52:48.770489",
    "2021-10-05 10:52:48.770492",
    "2021-10-05 10:52:48.770495"
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timeValues[$i]
}
